# Fruta / hortaliza, semanal
# The underlying data rows (2-21) were reshuffled: each row's Fecha (D),
# Calidad (L), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) were replaced with the
# values belonging to a different date/quality combination already present
# in the sheet. All other columns (A, B, C, E, F, G, H, I, J, K, Q, R, T)
# stay the same for every row. Below is the final target state for each
# data row, applied directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row : Fecha, Calidad, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg
$rows = @(
    @{ Row = 2;  D = 44532; L = "Primera"; M = 100; N = 10000; O = 10000; P = 10000; S = 5000 },
    @{ Row = 3;  D = 44532; L = "Segunda"; M = 100; N = 8000;  O = 8000;  P = 8000;  S = 4000 },
    @{ Row = 4;  D = 44609; L = "Primera"; M = 100; N = 6500;  O = 7000;  P = 6750;  S = 3375 },
    @{ Row = 5;  D = 44609; L = "Segunda"; M = 50;  N = 6000;  O = 6000;  P = 6000;  S = 3000 },
    @{ Row = 6;  D = 45014; L = "Primera"; M = 100; N = 7000;  O = 7500;  P = 7250;  S = 3625 },
    @{ Row = 7;  D = 44988; L = "Primera"; M = 200; N = 6000;  O = 7000;  P = 6500;  S = 3250 },
    @{ Row = 8;  D = 44195; L = "Primera"; M = 200; N = 3000;  O = 3500;  P = 3250;  S = 1625 },
    @{ Row = 9;  D = 44195; L = "Segunda"; M = 100; N = 2500;  O = 2500;  P = 2500;  S = 1250 },
    @{ Row = 10; D = 44216; L = "Primera"; M = 200; N = 3500;  O = 4000;  P = 3750;  S = 1875 },
    @{ Row = 11; D = 44216; L = "Segunda"; M = 100; N = 3000;  O = 3000;  P = 3000;  S = 1500 },
    @{ Row = 12; D = 44944; L = "Primera"; M = 100; N = 7000;  O = 8000;  P = 7500;  S = 3750 },
    @{ Row = 13; D = 44574; L = "Primera"; M = 200; N = 7000;  O = 8000;  P = 7500;  S = 3750 },
    @{ Row = 14; D = 44574; L = "Segunda"; M = 100; N = 6000;  O = 6000;  P = 6000;  S = 3000 },
    @{ Row = 15; D = 44910; L = "Primera"; M = 200; N = 7500;  O = 8000;  P = 7750;  S = 3875 },
    @{ Row = 16; D = 44602; L = "Primera"; M = 200; N = 6000;  O = 7000;  P = 6500;  S = 3250 },
    @{ Row = 17; D = 44602; L = "Segunda"; M = 100; N = 5000;  O = 5000;  P = 5000;  S = 2500 },
    @{ Row = 18; D = 44617; L = "Primera"; M = 200; N = 6000;  O = 7000;  P = 6500;  S = 3250 },
    @{ Row = 19; D = 44559; L = "Primera"; M = 200; N = 6000;  O = 7000;  P = 6500;  S = 3250 },
    @{ Row = 20; D = 44559; L = "Segunda"; M = 100; N = 5000;  O = 5000;  P = 5000;  S = 2500 },
    @{ Row = 21; D = 44917; L = "Primera"; M = 200; N = 7000;  O = 7500;  P = 7250;  S = 3625 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 4).Value  = $r.D   # D: Fecha
    $ws.Cells.Item($n, 12).Value = $r.L   # L: Calidad
    $ws.Cells.Item($n, 13).Value = $r.M   # M: Volumen
    $ws.Cells.Item($n, 14).Value = $r.N   # N: Precio minimo
    $ws.Cells.Item($n, 15).Value = $r.O   # O: Precio maximo
    $ws.Cells.Item($n, 16).Value = $r.P   # P: Precio promedio ponderado
    $ws.Cells.Item($n, 19).Value = $r.S   # S: Precio $/Kg
}
